$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 277.665479297809
$ws.Range("C2").Value = 188.806311936178
$ws.Range("D2").Value = 170.344052473418
$ws.Range("I2").Value = -70.3345207021912
$ws.Range("B3").Value = 248.288890625911
$ws.Range("C3").Value = 185.514477780561
$ws.Range("D3").Value = 167.99723675755
$ws.Range("I3").Value = 122.288890625911
$ws.Range("B4").Value = 192.996465743396
$ws.Range("I4").Value = 100.996465743396
$ws.Range("B5").Value = 210.465993195331
$ws.Range("I5").Value = 156.465993195331
$ws.Range("B6").Value = 238.028154326735
$ws.Range("I6").Value = 198.028154326735
$ws.Range("B7").Value = 262.552826401736
$ws.Range("I7").Value = 188.552826401736
$ws.Range("B8").Value = 272.354320567028
$ws.Range("I8").Value = 215.354320567028
$ws.Range("B9").Value = 238.906869506803
$ws.Range("I9").Value = 181.906869506803
$ws.Range("B10").Value = 226.3727454996
$ws.Range("I10").Value = 145.3727454996
$ws.Range("B11").Value = 207.029294210189
$ws.Range("I11").Value = 129.029294210189
$ws.Range("B12").Value = 224.006233758253
$ws.Range("I12").Value = 161.006233758253
$ws.Range("B13").Value = 200.604423664715
$ws.Range("I13").Value = 127.604423664715
$ws.Range("B14").Value = 263.027932310936
$ws.Range("I14").Value = 195.027932310936
$ws.Range("B15").Value = 235.813594676562
$ws.Range("I15").Value = 177.813594676562
$ws.Range("B16").Value = 190.028123387636
$ws.Range("I16").Value = 127.028123387636
$ws.Range("B17").Value = 206.103142438478
$ws.Range("I17").Value = 132.103142438478
$ws.Range("B18").Value = 231.162166792109
$ws.Range("I18").Value = 134.162166792109
$ws.Range("B19").Value = 253.331642609429
$ws.Range("I19").Value = 144.331642609429
$ws.Range("B20").Value = 261.026623422431
$ws.Range("I20").Value = 169.026623422431
$ws.Range("B21").Value = 230.429080373768
$ws.Range("I21").Value = 170.429080373768
$ws.Range("B22").Value = 219.113536962381
$ws.Range("I22").Value = 170.113536962381
$ws.Range("B23").Value = 202.712094157938
$ws.Range("I23").Value = 156.712094157938
$ws.Range("B24").Value = 218.115464100946
$ws.Range("I24").Value = 156.115464100946
$ws.Range("B25").Value = 197.452625333663
$ws.Range("I25").Value = 125.452625333663
$ws.Range("B26").Value = 254.692327031846
$ws.Range("I26").Value = 228.692327031846
$ws.Range("B27").Value = 229.43705795319
$ws.Range("I27").Value = 199.43705795319
$ws.Range("B28").Value = 187.133967645228
$ws.Range("I28").Value = 156.133967645228
$ws.Range("B29").Value = 203.349724010194
$ws.Range("I29").Value = 169.349724010194
$ws.Range("B30").Value = 226.876825960591
$ws.Range("I30").Value = 189.876825960591
$ws.Range("B31").Value = 246.345912501205
$ws.Range("I31").Value = 184.345912501205
$ws.Range("B32").Value = 252.09095850747
$ws.Range("I32").Value = 173.09095850747
$ws.Range("B33").Value = 223.631169703654
$ws.Range("I33").Value = 6.63116970365436
$ws.Range("B34").Value = 213.299097700413
$ws.Range("I34").Value = 117.299097700413
$ws.Range("B35").Value = 198.342132721394
$ws.Range("I35").Value = 106.342132721394
$ws.Range("B36").Value = 211.247570692797
$ws.Range("I36").Value = 137.247570692797
$ws.Range("B37").Value = 192.842713208345
$ws.Range("I37").Value = 113.842713208345
$ws.Range("B38").Value = 246.344460076058
$ws.Range("I38").Value = 97.3444600760585
$ws.Range("B39").Value = 222.465525832339
$ws.Range("I39").Value = 105.465525832339
$ws.Range("B40").Value = 183.334803586605
$ws.Range("I40").Value = 45.3348035866051
$ws.Range("B41").Value = 199.248821556037
$ws.Range("I41").Value = 10.2488215560367
$ws.Range("B42").Value = 221.881118962847
$ws.Range("I42").Value = 9.88111896284747
$ws.Range("B43").Value = 240.232480933666
$ws.Range("I43").Value = -23.7675190663343
$ws.Range("B44").Value = 245.082044099709
$ws.Range("I44").Value = -43.9179559002907
$ws.Range("B45").Value = 217.817374851696
$ws.Range("I45").Value = -16.1826251483036
$ws.Range("B46").Value = 207.947280003785
$ws.Range("I46").Value = 14.9472800037852
$ws.Range("B47").Value = 194.113984190323
$ws.Range("I47").Value = 1.11398419032309
$ws.Range("B48").Value = 206.419036461222
$ws.Range("I48").Value = 23.4190364612221
$ws.Range("B49").Value = 189.133675067448
$ws.Range("I49").Value = -55.8663249325515
